$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old)
    if ($found) {
        $rng.Text = $new
    } else {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "Starry Explorations: Unraveling the Cosmic Tapestry" "Shaping Tomorrow's Minds: The Profound Impact of History on Our Future"
Replace-Text "Lilyana Jireckova" "Eleanor Richards"
Replace-Text "l.jireckova@astroyale" "erichards@brightschools"
Replace-Text "With an insatiable curiosity about our place in the boundless universe, we embark on a captivating journey to uncover the secrets hidden amidst the shimmering stars and swirling galaxies" "In the vast tapestry of human existence, few disciplines hold such profound sway over our present and future as history"
Replace-Text " Our quest for knowledge has fueled our quest to unravel the cosmic tapestry, unlocking the mysteries that lie beyond Earth's atmosphere" " It serves as the beacon illuminating our path, shedding light on the past to illuminate the path towards a better tomorrow"
Replace-Text " Through brilliant minds, innovative technologies, and unwavering dedication, we strive to decode the intricate language of the cosmos, pushing the boundaries of our understanding and expanding the frontiers of human knowledge" " From the birth of great civilizations to the ebb and flow of empires, the study of history offers a panoramic view of human achievement and folly, imparting timeless lessons that shape our understanding of the world"
Replace-Text "Across civilizations and millennia, humanity has gazed upon the heavens with awe and wonder" "It is within the annals of history that we find the blueprint of our collective identity, a kaleidoscope of diverse cultures and traditions that have come together to create the rich mosaic of humanity"
Replace-Text " Ancient stargazers mapped constellations, divined celestial omens, and pondered the nature of the universe" " The study of these myriad threads that weave together the fabric of our existence fosters tolerance, understanding, and respect for the myriad ways in which civilizations and individuals have navigated the ebb and flow of time"
Replace-Text " Today, we continue this tradition, employing cutting-edge instruments and ingenious minds to explore the vast expanse of space. Our probes journey to distant worlds, uncovering clues to the solar system's formation and searching for signs of extraterrestrial life. Telescopes peer into the depths of space-time, revealing black holes, neutron stars, and galaxies billions of light-years away" " History thus becomes a potent force for unity, steering us towards a future where differences are embraced rather than feared"
Replace-Text "We seek to penetrate the veil of darkness and traverse the cosmos like never before" "Moreover, history provides an unparalleled lens through which we can examine the intricacies of human nature"
Replace-Text " With missions like the James Webb Space Telescope, we venture into uncharted territories, peering into the earliest moments of the universe and witnessing the birth of stars and galaxies" " As we delve into the lives of influential figures from across time, we gain insights into the motivations, passions, and fears that drive our actions"
Replace-Text " Through interplanetary missions, we search for habitable environments and scrutinize distant exoplanets, hoping to discover worlds beyond our own" " This introspective journey helps us better comprehend our own place in the grand scheme of things, instilling empathy and compassion for the human condition"
Replace-Text " Every new discovery, every step forward, brings us closer to comprehending the grand symphony of the universe, revealing its infinite beauty and complexity" " Through this prism, history guides us towards becoming more thoughtful decision-makers, whose actions are informed by a deep appreciation for both the complexity of the past and the potential of the future"
Replace-Text "Our cosmic voyage leads us through celestial wonders, pushing the boundaries of scientific understanding" "In this essay, we have explored the profound impact of history on our present and future"
Replace-Text " We witness the birth of stars, explore exoplanets, and uncover the mysteries of black holes" " We have delved into the role of history in shaping our understanding of the world, inculcating tolerance and respect for diverse cultures, and providing insights into the complexities of human nature"
Replace-Text " Immersed in the language of the universe, we strive for a deeper connection with the cosmos, searching for life beyond Earth and yearning to unravel the ultimate secrets of " " History, therefore, serves as a vital tool for creating more thoughtful and compassionate "
Replace-Text "existence. With each discovery, we are awestruck by the universe's limitless majesty, humbled by its grandeur, and inspired to continue our relentless pursuit of knowledge" "citizens who are equipped with the knowledge and skills necessary to forge a better future for all"

Write-Host "DONE"
